$wb = $excel.ActiveWorkbook

$wsAll = $wb.Worksheets.Item("All Tests")
$wsSprint = $wb.Worksheets.Item("Sprint 1")

# --- Content edits -------------------------------------------------------
# "invite functionality" row (All Tests, row 10): tweak wording of the
# acceptance-criteria text for backlog item #13 (invite/join flow).
$wsAll.Range("C10").Value2 = "a room owner should be able to invite people by invite link or username, users should be able to request to join a room by entering the room code or searching for rooms, the room owner should be able to accept or deny the request and the user should be able to see pending requests"

# "recipe functionality" row (All Tests, row 6): tweak wording of the
# acceptance-criteria text for backlog item #14.
$wsAll.Range("C6").Value2 = "satisfy criteria of #14, when creating a meal, a user should be able to choose from saved or public recipes or create a new recipe, a recipe should contain a list of ingredients, steps of preparation (description) and time to prepare, a recipe should be either private or public, a recipe should be able to be shared per link"

# --- Selection / active-sheet state --------------------------------------
# The workbook now opens on "All Tests" (was "Sprint 1"), with a new
# selected cell on each sheet.
$wsSprint.Range("C17").Select()
$wsAll.Activate()
$wsAll.Range("C19").Select()

# --- Column sizing ---------------------------------------------------------
# Columns were widened (auto-fit) to accommodate the revised text.
$wsAll.Columns.Item(2).AutoFit()
$wsAll.Columns.Item(3).AutoFit()
$wsSprint.Columns.Item(2).AutoFit()
$wsSprint.Columns.Item(3).AutoFit()

Write-Host "edit applied"
